# Auto-generated Excel COM-interop script to apply the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '56.615.43'
Set-TextCell 'E2' '  +2.94%  '
Set-TextCell 'D3' '2.322.68'
Set-TextCell 'E3' '  +1.03%  '
Set-TextCell 'D5' '520.14'
Set-TextCell 'E5' '  +2.41%  '
Set-TextCell 'D6' '134.74'
Set-TextCell 'E6' '  +3.51%  '
Set-TextCell 'D7' '0.995'
Set-TextCell 'E7' '  +0.09%  '
Set-TextCell 'D8' '0.537'
Set-TextCell 'E8' '  +0.96%  '
Set-TextCell 'D9' '2.347.80'
Set-TextCell 'E9' '  +1.04%  '
Set-TextCell 'E10' '  +5.77%  '
Set-TextCell 'E11' '  -0.84%  '
Set-TextCell 'D12' '5.26'
Set-TextCell 'E12' '  +3.86%  '
Set-TextCell 'E13' '  +0.01%  '
Set-TextCell 'D14' '23.80'
Set-TextCell 'E14' '  -0.95%  '
Set-TextCell 'D15' '2.741.32'
Set-TextCell 'E15' '  +1.16%  '
Set-TextCell 'D16' '56.675.86'
Set-TextCell 'E16' '  +3.17%  '
Set-TextCell 'E17' '  +2.01%  '
Set-TextCell 'D18' '2.328.17'
Set-TextCell 'E18' '  +1.64%  '
Set-TextCell 'D19' '10.45'
Set-TextCell 'E19' '  -2.38%  '
Set-TextCell 'D20' '4.22'
Set-TextCell 'E20' '  +0.70%  '
Set-TextCell 'E21' '  +3.71%  '
Set-TextCell 'D22' '6.55'
Set-TextCell 'E22' '  -2.06%  '
Set-TextCell 'D23' '1.00'
Set-TextCell 'E23' '  +0.38%  '
Set-TextCell 'D24' '60.57'
Set-TextCell 'E24' '  +0.07%  '
Set-TextCell 'D25' '0.163'
Set-TextCell 'E25' '  +7.24%  '
Set-TextCell 'D26' '0.994'
Set-TextCell 'E26' '  +0.14%  '
Set-TextCell 'D27' '7.88'
Set-TextCell 'E27' '  +4.39%  '
Set-TextCell 'D28' '1.27'
Set-TextCell 'E28' '  +10.02%  '
Set-TextCell 'D29' '0.0₃0749'
Set-TextCell 'E29' '  +5.12%  '
Set-TextCell 'D30' '170.07'
Set-TextCell 'E30' '  -1.70%  '
Set-TextCell 'D31' '1.73'
Set-TextCell 'E31' '  +5.45%  '
Set-TextCell 'D32' '6.18'
Set-TextCell 'E32' '  +0.14%  '
Set-TextCell 'D33' '18.25'
Set-TextCell 'E33' '  +0.82%  '
Set-TextCell 'E34' '  +0.05%  '
Set-TextCell 'E35' '  -0.09%  '
Set-TextCell 'D36' '1.25'
Set-TextCell 'E36' '  +1.11%  '
Set-TextCell 'D37' '0.922'
Set-TextCell 'E37' '  +0.45%  '
Set-TextCell 'D38' '4.01'
Set-TextCell 'E38' '  +2.75%  '
Set-TextCell 'E39' '  +7.33%  '
Set-TextCell 'D40' '37.96'
Set-TextCell 'E40' '  +3.12%  '
Set-TextCell 'E41' '  +0.04%  '
Set-TextCell 'D42' '3.58'
Set-TextCell 'E42' '  +4.28%  '
Set-TextCell 'D43' '137.26'
Set-TextCell 'E43' '  +0.96%  '
Set-TextCell 'B44' 'RenderToken'
Set-TextCell 'C44' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D44' '5.22'
Set-TextCell 'B45' 'Bittensor'
Set-TextCell 'C45' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D45' '277.74'
Set-TextCell 'E45' '  +6.07%  '
Set-TextCell 'D46' '0.0933'
Set-TextCell 'E46' '  +2.08%  '
Set-TextCell 'D47' '0.0503'
Set-TextCell 'E47' '  -0.28%  '
Set-TextCell 'E48' '  +1.10%  '
Set-TextCell 'E49' '  +3.38%  '
Set-TextCell 'B50' 'InjectiveProtocol'
Set-TextCell 'C50' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D50' '17.90'
Set-TextCell 'E50' '  +7.81%  '
Set-TextCell 'B51' 'Polygon'
Set-TextCell 'C51' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D51' '0.379'
Set-TextCell 'E51' '  +0.12%  '
